$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 928.44446
$ws.Range("I58").Value = 397
$ws.Range("J58").Value = 1459.8889
$ws.Range("K58").Value = 1191
$ws.Range("L58").Value = 4379.6667
$ws.Range("M58").Value = -1041
$ws.Range("N58").Value = -4679.6667

$ws.Range("H69").Value = 1737.25
$ws.Range("I69").Value = 1737.25
$ws.Range("K69").Value = 5211.75
$ws.Range("M69").Value = -4337.75

$ws.Range("H72").Value = 1737.25
$ws.Range("I72").Value = 1737.25
$ws.Range("K72").Value = 15635.25
$ws.Range("M72").Value = -11267.25

$ws.Range("H74").Value = 5120
$ws.Range("J74").Value = 3900
$ws.Range("L74").Value = 3900
$ws.Range("N74").Value = -5772

$ws.Range("H76").Value = 1305077.4
$ws.Range("I76").Value = 2131942.5
$ws.Range("K76").Value = 2131942.5
$ws.Range("M76").Value = -2131627.5

$ws.Range("H77").Value = 5120
$ws.Range("J77").Value = 3900
$ws.Range("L77").Value = 19500
$ws.Range("N77").Value = -28860

$ws.Range("H79").Value = 1305077.4
$ws.Range("I79").Value = 2131942.5
$ws.Range("K79").Value = 2131942.5
$ws.Range("M79").Value = -2130850.5

$ws.Range("H116").Value = 23560
$ws.Range("I116").Value = 35933.332
$ws.Range("J116").Value = 5000
$ws.Range("K116").Value = 35933.332
$ws.Range("L116").Value = 5000
$ws.Range("M116").Value = -32491.332
$ws.Range("N116").Value = -11884

$ws.Range("H132").Value = 1282.7407
$ws.Range("I132").Value = 1172.25
$ws.Range("K132").Value = 3516.75
$ws.Range("M132").Value = -986.75

$ws.Range("H138").Value = 2991.08
$ws.Range("J138").Value = 2828.5715
$ws.Range("L138").Value = 8485.7145
$ws.Range("N138").Value = -18765.7145

$ws.Range("H141").Value = 1650655.2
$ws.Range("I141").Value = 2547914
$ws.Range("K141").Value = 7643742
$ws.Range("M141").Value = -7638562

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5085.684
$ws.Range("I32").Value = 3702.4
$ws.Range("J32").Value = 10273
$ws.Range("K32").Value = 3702.4
$ws.Range("L32").Value = 10273
$ws.Range("M32").Value = -3415.4
$ws.Range("N32").Value = -10847

$ws.Range("M37").ClearContents()
$ws.Range("H37").Value = 14028.5
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 14028.5
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 14028.5
$ws.Range("N37").Value = -14574.5

$ws.Range("H74").Value = 1208.591
$ws.Range("I74").Value = 793.94446
$ws.Range("J74").Value = 3074.5
$ws.Range("K74").Value = 793.94446
$ws.Range("L74").Value = 3074.5
$ws.Range("M74").Value = 80.05553999999995
$ws.Range("N74").Value = -4822.5

$ws.Range("H77").Value = 1208.591
$ws.Range("I77").Value = 793.94446
$ws.Range("J77").Value = 3074.5
$ws.Range("K77").Value = 3969.7223
$ws.Range("L77").Value = 15372.5
$ws.Range("M77").Value = 398.2776999999996
$ws.Range("N77").Value = -24108.5

$ws.Range("H110").Value = 2454.3635
$ws.Range("I110").Value = 1434
$ws.Range("J110").Value = 5175.3335
$ws.Range("K110").Value = 1434
$ws.Range("L110").Value = 5175.3335
$ws.Range("M110").Value = 611
$ws.Range("N110").Value = -9265.333500000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2252.6667
$ws.Range("I20").Value = 2315.7
$ws.Range("J20").Value = 2126.6
$ws.Range("K20").Value = 2315.7
$ws.Range("L20").Value = 2126.6
$ws.Range("M20").Value = -2068.7
$ws.Range("N20").Value = -2620.6

$ws.Range("H35").Value = 35000
$ws.Range("J35").Value = 35000
$ws.Range("L35").Value = 35000
$ws.Range("N35").Value = -35620

$ws.Range("H94").Value = 335.96667
$ws.Range("I94").Value = 342.3793
$ws.Range("K94").Value = 342.3793
$ws.Range("M94").Value = 108.6207

$ws.Range("H99").Value = 1543.2307
$ws.Range("I99").Value = 1306.2
$ws.Range("K99").Value = 1306.2
$ws.Range("M99").Value = 191.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 748.1
$ws.Range("I16").Value = 735.375
$ws.Range("K16").Value = 735.375
$ws.Range("M16").Value = -448.375

$ws.Range("N55").ClearContents()
$ws.Range("H55").Value = 9000
$ws.Range("I55").Value = 9000
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 9000
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -8685

$ws.Range("H58").Value = 2559103.2
$ws.Range("I58").Value = 3624445.2
$ws.Range("K58").Value = 3624445.2
$ws.Range("M58").Value = -3624242.2

$ws.Range("H113").Value = 748.1
$ws.Range("I113").Value = 735.375
$ws.Range("K113").Value = 735.375
$ws.Range("M113").Value = 1434.625

$ws.Range("H132").Value = 2190.375
$ws.Range("I132").Value = 1172.375
$ws.Range("K132").Value = 3517.125
$ws.Range("M132").Value = -987.125

$ws.Range("H136").Value = 2559103.2
$ws.Range("I136").Value = 3624445.2
$ws.Range("K136").Value = 10873335.6
$ws.Range("M136").Value = -10870785.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 930
$ws.Range("J5").Value = 936
$ws.Range("L5").Value = 2808
$ws.Range("N5").Value = -3032

$ws.Range("H87").Value = 11569.143
$ws.Range("I87").Value = 6196.8
$ws.Range("K87").Value = 18590.4
$ws.Range("M87").Value = -17342.4

$ws.Range("H90").Value = 11569.143
$ws.Range("I90").Value = 6196.8
$ws.Range("K90").Value = 55771.2
$ws.Range("M90").Value = -49531.2

$ws.Range("H135").Value = 930
$ws.Range("J135").Value = 936
$ws.Range("L135").Value = 8424
$ws.Range("N135").Value = -13494

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 12971.111
$ws.Range("I70").Value = 19928.2
$ws.Range("J70").Value = 4274.75
$ws.Range("K70").Value = 19928.2
$ws.Range("L70").Value = 4274.75
$ws.Range("M70").Value = -19658.2
$ws.Range("N70").Value = -4814.75

$ws.Range("H73").Value = 12971.111
$ws.Range("I73").Value = 19928.2
$ws.Range("J73").Value = 4274.75
$ws.Range("K73").Value = 19928.2
$ws.Range("L73").Value = 4274.75
$ws.Range("M73").Value = -18992.2
$ws.Range("N73").Value = -6146.75

$ws.Range("H80").Value = 2649
$ws.Range("I80").Value = 2448.25
$ws.Range("K80").Value = 2448.25
$ws.Range("M80").Value = -1450.25

$ws.Range("H83").Value = 2649
$ws.Range("I83").Value = 2448.25
$ws.Range("K83").Value = 12241.25
$ws.Range("M83").Value = -7249.25

$ws.Range("H132").Value = 2407931.8
$ws.Range("I132").Value = 2750422
$ws.Range("K132").Value = 8251266
$ws.Range("M132").Value = -8248736

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3116
$ws.Range("J22").Value = 2099
$ws.Range("L22").Value = 2099
$ws.Range("N22").Value = -2689

$ws.Range("H27").Value = 3116
$ws.Range("J27").Value = 2099
$ws.Range("L27").Value = 2099
$ws.Range("N27").Value = -2313

$ws.Range("H46").Value = 1403.0834
$ws.Range("J46").Value = 1548.7
$ws.Range("L46").Value = 1548.7
$ws.Range("N46").Value = -1924.7

$ws.Range("H132").Value = 2637.3572
$ws.Range("I132").Value = 2262.7778
$ws.Range("K132").Value = 6788.3334
$ws.Range("M132").Value = -4258.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2998.5
$ws.Range("I62").Value = 2998
$ws.Range("K62").Value = 2998
$ws.Range("M62").Value = -2374

$ws.Range("H65").Value = 2998.5
$ws.Range("I65").Value = 2998
$ws.Range("K65").Value = 14990
$ws.Range("L65").Value = 14995
$ws.Range("M65").Value = -11870

$ws.Range("H126").Value = 12388.846
$ws.Range("J126").Value = 6556.7144
$ws.Range("L126").Value = 19670.1432
$ws.Range("N126").Value = -24610.1432

$ws.Range("H132").Value = 4099.3
$ws.Range("J132").Value = 5199.2
$ws.Range("L132").Value = 15597.6
$ws.Range("N132").Value = -20657.6

$ws.Range("H136").Value = 17923646
$ws.Range("I136").Value = 30866680
$ws.Range("K136").Value = 92600040
$ws.Range("M136").Value = -92597490
